$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the color image path and main image path for the product row,
# replacing the old broken absolute paths with the new ones.
$ws.Range("L2").Value = "D:\imgDATN\2.png"
$ws.Range("M2").Value = "D:\imgDATN\3.png"

# Reflect the cell that was selected/active when the edit was made.
[void]$ws.Range("N4").Select()
